# feat: add 2022-Q4 data
#
# Before: sheets = [ "总计", "2022-Q1" ]
# After:  sheets = [ "总计", "2022-Q4", "2022-Q1" ]
#   - "总计" gains a row for the new 2022-Q4 totals (existing 2022-Q1 row shifts down).
#   - A new "2022-Q4" sheet is inserted holding the new quarter's fund holdings.
#   - The original "2022-Q1" sheet (and its data) is preserved unchanged, just moved
#     one slot to the right.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q1Sheet = $wb.Worksheets.Item("2022-Q1")

# ---------------------------------------------------------------------------
# 1) Clone the existing "2022-Q1" sheet and park the clone at the end of the
#    workbook. That clone keeps the old fund-holdings data/formatting intact
#    and will become the new "2022-Q1" tab (after "2022-Q4").
# ---------------------------------------------------------------------------
$q1Sheet.Copy($null, $q1Sheet)
$clone = $wb.Worksheets.Item($wb.Worksheets.Count)
$clone.Name = "2022-Q1 NEW"

# ---------------------------------------------------------------------------
# 2) Repurpose the original "2022-Q1" sheet (still sitting in slot 2, right
#    after "总计") into the new "2022-Q4" sheet: rename it and replace its
#    contents with the new quarter's fund-holdings table.
# ---------------------------------------------------------------------------
$q1Sheet.Name = "2022-Q4"
$q4Sheet = $q1Sheet

$q4Sheet.Cells.Clear()

# Match the page margins used by the "总计" sheet / the rest of the workbook.
$q4Sheet.PageSetup.LeftMargin = 54
$q4Sheet.PageSetup.RightMargin = 54
$q4Sheet.PageSetup.TopMargin = 72
$q4Sheet.PageSetup.BottomMargin = 72
$q4Sheet.PageSetup.HeaderMargin = 36
$q4Sheet.PageSetup.FooterMargin = 36

# Pull the header / index-column styling from "总计" (style used on its
# B1:D1 header row and its A2 index cell) onto the new sheet.
$totalSheet.Range("B1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)
$totalSheet.Range("A2").Copy()
$q4Sheet.Range("A2:A5").PasteSpecial(-4122)

$q4Sheet.Range("B1").Value = "基金代码"
$q4Sheet.Range("C1").Value = "基金名称"
$q4Sheet.Range("D1").Value = "基金规模"
$q4Sheet.Range("E1").Value = "股票总仓位"
$q4Sheet.Range("F1").Value = "仓位占比"
$q4Sheet.Range("G1").Value = "持有市值(亿元)"
$q4Sheet.Range("H1").Value = "仓位排名"

$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("B2").Value = "'008638"
$q4Sheet.Range("C2").Value = "广发科技创新混合A"
$q4Sheet.Range("D2").Value = "'33.72"
$q4Sheet.Range("E2").Value = "'91.92"
$q4Sheet.Range("F2").Value = "'6.27"
$q4Sheet.Range("G2").Value = "'2.1142"
$q4Sheet.Range("H2").Value = 4

$q4Sheet.Range("A3").Value = 1
$q4Sheet.Range("B3").Value = "'013533"
$q4Sheet.Range("C3").Value = "广发科技创新混合C"
$q4Sheet.Range("D3").Value = "'2.90"
$q4Sheet.Range("E3").Value = "'91.92"
$q4Sheet.Range("F3").Value = "'6.27"
$q4Sheet.Range("G3").Value = "'0.1818"
$q4Sheet.Range("H3").Value = 4

$q4Sheet.Range("A4").Value = 2
$q4Sheet.Range("B4").Value = "'006429"
$q4Sheet.Range("C4").Value = "诺安恒鑫混合"
$q4Sheet.Range("D4").Value = "'0.92"
$q4Sheet.Range("E4").Value = "'66.36"
$q4Sheet.Range("F4").Value = "'4.26"
$q4Sheet.Range("G4").Value = "'0.0392"
$q4Sheet.Range("H4").Value = 7

$q4Sheet.Range("A5").Value = 3
$q4Sheet.Range("B5").Value = "'620002"
$q4Sheet.Range("C5").Value = "金元顺安成长动力混合"
$q4Sheet.Range("D5").Value = "'0.35"
$q4Sheet.Range("E5").Value = "'72.02"
$q4Sheet.Range("F5").Value = "'3.89"
$q4Sheet.Range("G5").Value = "'0.0136"
$q4Sheet.Range("H5").Value = 2

# ---------------------------------------------------------------------------
# 3) Rename the clone back to "2022-Q1" - it already carries forward the
#    original fund-holdings data untouched.
# ---------------------------------------------------------------------------
$clone.Name = "2022-Q1"

# ---------------------------------------------------------------------------
# 4) Update the "总计" summary sheet: the old row 2 (2022-Q1 totals) shifts
#    down to row 3, and row 2 now holds the new 2022-Q4 totals.
# ---------------------------------------------------------------------------
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q1"
$totalSheet.Range("C3").Value = 3
$totalSheet.Range("D3").Value = 0.07000000000000001

$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 2.35
